# Applies the changes described by the commit "fix: slide changes; missing mp4":
#   1. Update the fixed Header & Footer date placeholder text from 6/9/24 to
#      6/10/24 on the slide master and on every slide layout.
#   2. Slide 2: "docker-compose up" -> "docker compose up" (en dash kept)
#   3. Slide 4: "They have supporting, short videos in Dropbox." ->
#      "They have supporting, short videos in GitHub."

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Date placeholder (6/9/24 -> 6/10/24) on the slide master + all layouts.
# ---------------------------------------------------------------------------
$oldDate = "6/9/24"
$newDate = "6/10/24"

$master = $p.SlideMaster

for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.HasTextFrame -eq -1) {
        if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.HasTextFrame -eq -1) {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# ---------------------------------------------------------------------------
# 2. Slide 2: "docker-compose up -d" -> "docker compose up -d" (en dash kept)
# ---------------------------------------------------------------------------
$slide2 = $p.Slides.Item(2)
$shape2 = $slide2.Shapes.Item(2)
$tr2 = $shape2.TextFrame.TextRange
$run = $tr2.Find("docker-compose up –d", 0)
$run.Text = "docker compose up –d"

# ---------------------------------------------------------------------------
# 3. Slide 4: "... videos in Dropbox." -> "... videos in GitHub."
# ---------------------------------------------------------------------------
$slide4 = $p.Slides.Item(4)
$shape4 = $slide4.Shapes.Item(2)
$tr4 = $shape4.TextFrame.TextRange
$found = $tr4.Find("in Dropbox.", 0)
$found.Text = "in GitHub."
